# Auto-generated: update 2025 (column L) violent-crime figures for 2025-06-03
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Cells.Item(2, 12).Value = 2668
$ws.Cells.Item(3, 12).Value = 2693
$ws.Cells.Item(4, 12).Value = 715
$ws.Cells.Item(5, 12).Value = 158
$ws.Cells.Item(6, 12).Value = 2413
$ws.Cells.Item(7, 12).Value = 8647

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Cells.Item(5, 12).Value = 34
$ws.Cells.Item(6, 12).Value = 66
$ws.Cells.Item(7, 12).Value = 287
$ws.Cells.Item(8, 12).Value = 550
$ws.Cells.Item(9, 12).Value = 56
$ws.Cells.Item(10, 12).Value = 57
$ws.Cells.Item(11, 12).Value = 153
$ws.Cells.Item(15, 12).Value = 61
$ws.Cells.Item(19, 12).Value = 242
$ws.Cells.Item(23, 12).Value = 86
$ws.Cells.Item(27, 12).Value = 86
$ws.Cells.Item(29, 12).Value = 458
$ws.Cells.Item(31, 12).Value = 82
$ws.Cells.Item(33, 12).Value = 388
$ws.Cells.Item(37, 12).Value = 318
$ws.Cells.Item(42, 12).Value = 284
$ws.Cells.Item(44, 12).Value = 64
$ws.Cells.Item(45, 12).Value = 16
$ws.Cells.Item(47, 12).Value = 65
$ws.Cells.Item(49, 12).Value = 48
$ws.Cells.Item(51, 12).Value = 104
$ws.Cells.Item(52, 12).Value = 175
$ws.Cells.Item(54, 12).Value = 172
$ws.Cells.Item(55, 12).Value = 83
$ws.Cells.Item(60, 12).Value = 51
$ws.Cells.Item(63, 12).Value = 26
$ws.Cells.Item(65, 12).Value = 160
$ws.Cells.Item(66, 12).Value = 20
$ws.Cells.Item(67, 12).Value = 318
$ws.Cells.Item(68, 12).Value = 27
$ws.Cells.Item(72, 12).Value = 36
$ws.Cells.Item(73, 12).Value = 71
$ws.Cells.Item(77, 12).Value = 51
$ws.Cells.Item(79, 12).Value = 231
$ws.Cells.Item(83, 12).Value = 203
$ws.Cells.Item(85, 12).Value = 449
$ws.Cells.Item(86, 12).Value = 62
$ws.Cells.Item(87, 12).Value = 26
$ws.Cells.Item(88, 12).Value = 111
$ws.Cells.Item(89, 12).Value = 112
$ws.Cells.Item(91, 12).Value = 122
$ws.Cells.Item(94, 12).Value = 105
$ws.Cells.Item(95, 12).Value = 119
$ws.Cells.Item(96, 12).Value = 84
$ws.Cells.Item(99, 12).Value = 143
$ws.Cells.Item(101, 12).Value = 8647

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Cells.Item(2, 12).Value = 32
$ws.Cells.Item(7, 12).Value = 84

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Cells.Item(2, 12).Value = 85
$ws.Cells.Item(3, 12).Value = 92
$ws.Cells.Item(7, 12).Value = 287

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Cells.Item(3, 12).Value = 48
$ws.Cells.Item(4, 12).Value = 12
$ws.Cells.Item(6, 12).Value = 37
$ws.Cells.Item(7, 12).Value = 153

$ws = $wb.Worksheets.Item("Uptown")
$ws.Cells.Item(3, 12).Value = 29
$ws.Cells.Item(7, 12).Value = 112

$ws = $wb.Worksheets.Item("South Shore")
$ws.Cells.Item(4, 12).Value = 36
$ws.Cells.Item(6, 12).Value = 91
$ws.Cells.Item(7, 12).Value = 449

$ws = $wb.Worksheets.Item("Little Village")
$ws.Cells.Item(3, 12).Value = 51
$ws.Cells.Item(6, 12).Value = 48
$ws.Cells.Item(7, 12).Value = 175

$ws = $wb.Worksheets.Item("Austin")
$ws.Cells.Item(2, 12).Value = 156
$ws.Cells.Item(3, 12).Value = 182
$ws.Cells.Item(4, 12).Value = 37
$ws.Cells.Item(6, 12).Value = 155
$ws.Cells.Item(7, 12).Value = 550

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Cells.Item(2, 12).Value = 63
$ws.Cells.Item(3, 12).Value = 87
$ws.Cells.Item(7, 12).Value = 203

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Cells.Item(3, 12).Value = 124
$ws.Cells.Item(7, 12).Value = 388

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Cells.Item(6, 12).Value = 23
$ws.Cells.Item(7, 12).Value = 119

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Cells.Item(6, 12).Value = 100
$ws.Cells.Item(7, 12).Value = 318

$ws = $wb.Worksheets.Item("New City")
$ws.Cells.Item(3, 12).Value = 50
$ws.Cells.Item(7, 12).Value = 160

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Cells.Item(3, 12).Value = 61
$ws.Cells.Item(6, 12).Value = 34
$ws.Cells.Item(7, 12).Value = 143

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Cells.Item(3, 12).Value = 20
$ws.Cells.Item(7, 12).Value = 82

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Cells.Item(3, 12).Value = 113
$ws.Cells.Item(7, 12).Value = 318

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Cells.Item(3, 12).Value = 7
$ws.Cells.Item(6, 12).Value = 20
$ws.Cells.Item(7, 12).Value = 48

$ws = $wb.Worksheets.Item("Loop")
$ws.Cells.Item(2, 12).Value = 38
$ws.Cells.Item(6, 12).Value = 89
$ws.Cells.Item(7, 12).Value = 172

$ws = $wb.Worksheets.Item("Englewood")
$ws.Cells.Item(2, 12).Value = 147
$ws.Cells.Item(3, 12).Value = 168
$ws.Cells.Item(7, 12).Value = 458

$ws = $wb.Worksheets.Item("Chatham")
$ws.Cells.Item(3, 12).Value = 74
$ws.Cells.Item(6, 12).Value = 76
$ws.Cells.Item(7, 12).Value = 242

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Cells.Item(2, 12).Value = 27
$ws.Cells.Item(7, 12).Value = 64

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Cells.Item(2, 12).Value = 28
$ws.Cells.Item(3, 12).Value = 19
$ws.Cells.Item(7, 12).Value = 66

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Cells.Item(2, 12).Value = 79
$ws.Cells.Item(3, 12).Value = 86
$ws.Cells.Item(6, 12).Value = 85
$ws.Cells.Item(7, 12).Value = 284

$ws = $wb.Worksheets.Item("Avondale")
$ws.Cells.Item(2, 12).Value = 27
$ws.Cells.Item(7, 12).Value = 57

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Cells.Item(6, 12).Value = 19
$ws.Cells.Item(7, 12).Value = 83

$ws = $wb.Worksheets.Item("Douglas")
$ws.Cells.Item(6, 12).Value = 17
$ws.Cells.Item(7, 12).Value = 86

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Cells.Item(2, 12).Value = 49
$ws.Cells.Item(3, 12).Value = 45
$ws.Cells.Item(6, 12).Value = 16
$ws.Cells.Item(7, 12).Value = 122

$ws = $wb.Worksheets.Item("Roseland")
$ws.Cells.Item(2, 12).Value = 72
$ws.Cells.Item(7, 12).Value = 231

$ws = $wb.Worksheets.Item("West Loop")
$ws.Cells.Item(5, 12).Value = 3
$ws.Cells.Item(6, 12).Value = 34
$ws.Cells.Item(7, 12).Value = 105

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Cells.Item(6, 12).Value = 16
$ws.Cells.Item(7, 12).Value = 65

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Cells.Item(4, 12).Value = 6
$ws.Cells.Item(7, 12).Value = 61

$ws = $wb.Worksheets.Item("North Center")
$ws.Cells.Item(4, 12).Value = 5
$ws.Cells.Item(7, 12).Value = 20

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Cells.Item(2, 12).Value = 15
$ws.Cells.Item(7, 12).Value = 56

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Cells.Item(2, 12).Value = 31
$ws.Cells.Item(7, 12).Value = 71

$ws = $wb.Worksheets.Item("United Center")
$ws.Cells.Item(2, 12).Value = 31
$ws.Cells.Item(7, 12).Value = 111

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Cells.Item(6, 12).Value = 15
$ws.Cells.Item(7, 12).Value = 34

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Cells.Item(6, 12).Value = 23
$ws.Cells.Item(7, 12).Value = 86

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Cells.Item(4, 12).Value = 36
$ws.Cells.Item(7, 12).Value = 62

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Cells.Item(3, 12).Value = 30
$ws.Cells.Item(4, 12).Value = 16
$ws.Cells.Item(7, 12).Value = 104

$ws = $wb.Worksheets.Item("North Park")
$ws.Cells.Item(2, 12).Value = 6
$ws.Cells.Item(7, 12).Value = 27

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Cells.Item(3, 12).Value = 18
$ws.Cells.Item(7, 12).Value = 51

$ws = $wb.Worksheets.Item("Old Town")
$ws.Cells.Item(6, 12).Value = 9
$ws.Cells.Item(7, 12).Value = 36

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Cells.Item(5, 12).Value = 2
$ws.Cells.Item(7, 12).Value = 51

$ws = $wb.Worksheets.Item("Jackson Park")
$ws.Cells.Item(6, 12).Value = 5
$ws.Cells.Item(7, 12).Value = 16

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Cells.Item(3, 12).Value = 3
$ws.Cells.Item(7, 12).Value = 26
